$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new "profile refresh" test case rows (20-22), mirroring the
# existing layout: B=Platform, C=Method Name, D=Role, E=Key, F=Value
$ws.Range("B20").Value = "Web"
$ws.Range("C20").Value = "verifyPrimeClassOnProfileRefresh"
$ws.Range("D20").Value = "Student"
$ws.Range("E20").Value = "Prime Classes"
$ws.Range("F20").Value = "Class 11 Class 12"

$ws.Range("B21").Value = "Android"
$ws.Range("C21").Value = "verifyPrimeClassOnProfileRefresh"
$ws.Range("D21").Value = "Student"
$ws.Range("E21").Value = "Prime Classes"
$ws.Range("F21").Value = "Class 11 Class 12"

$ws.Range("B22").Value = "iOS"
$ws.Range("C22").Value = "verifyPrimeClassOnProfileRefresh"
$ws.Range("D22").Value = "Student"
$ws.Range("E22").Value = "Prime Classes"
$ws.Range("F22").Value = "Class 11 Class 12"

# Column C needs to widen to fit the new, longer method name
$ws.Columns("C:C").ColumnWidth = 31 + 1/6

# Move the selection/view down to the newly added last cell
[void]$ws.Range("F22").Select()
